$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (row height/format untouched) ---
# Raw OOXML widths 35/34/16/15/31 correspond to these COM ColumnWidth values
# (engine applies a fixed +0.83 char padding when round-tripping through xlsx).
$ws.Columns.Item(3).ColumnWidth = 34.17   # C: 90 -> 35
$ws.Columns.Item(4).ColumnWidth = 33.17   # D: 56 -> 34
$ws.Columns.Item(6).ColumnWidth = 15.17   # F: 17 -> 16
$ws.Columns.Item(7).ColumnWidth = 14.17   # G: 16 -> 15
$ws.Columns.Item(8).ColumnWidth = 30.17   # H: 57 -> 31

# --- Row data: replace rows 2-10 with the new scrape, append rows 11-13 ---
$rows = @(
    @(2, "1328069", "https://aiesec.org/opportunity/global-talent/1328069", "Food & Beverage - Intern", "Nugegoda, Sri Lanka", "No", "2 applicants", "3 - 6 Months", "Raffles Consolidated Pvt Ltd"),
    @(3, "1328052", "https://aiesec.org/opportunity/global-talent/1328052", "Digital Marketing Intern", "Sfax, Tunisie", "No", "1 applicant", "9 - 12 Weeks", "WIT center"),
    @(4, "1328041", "https://aiesec.org/opportunity/global-talent/1328041", "IT/ ERP System Developers (Odoo)", "Nugegoda, Sri Lanka", "No", "3 applicants", "3 - 6 Months", "Epigro Pvt Ltd"),
    @(5, "1328023", "https://aiesec.org/opportunity/global-talent/1328023", "Marketing - Intern", "Nugegoda, Sri Lanka", "No", "3 applicants", "3 - 6 Months", "Raffles Consolidated Pvt Ltd"),
    @(6, "1328019", "https://aiesec.org/opportunity/global-talent/1328019", "Kitchen Operations - Intern", "Nugegoda, Sri Lanka", "No", "1 applicant", "3 - 6 Months", "Raffles Consolidated Pvt Ltd"),
    @(7, "1328015", "https://aiesec.org/opportunity/global-talent/1328015", "Cybersecurity Intern", "Phagwara, Punjab, India", "No", "1 applicant", "3 - 6 Months", "GNA University"),
    @(8, "1328013", "https://aiesec.org/opportunity/global-talent/1328013", "Machine Learning Intern", "Phagwara, Punjab, India", "No", "1 applicant", "3 - 6 Months", "GNA University"),
    @(9, "1327992", "https://aiesec.org/opportunity/global-talent/1327992", "Cloud Engineering Intern", "Phagwara, Punjab, India", "No", "0 applicants", "3 - 6 Months", "GNA University"),
    @(10, "1327958", "https://aiesec.org/opportunity/global-talent/1327958", "Graphic designer", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "3 - 6 Months", "The Paddock"),
    @(11, "1327957", "https://aiesec.org/opportunity/global-talent/1327957", "Account Manager", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "3 - 6 Months", "The Paddock"),
    @(12, "1327951", "https://aiesec.org/opportunity/global-talent/1327951", "Account Executive", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Gipfel"),
    @(13, "1327149", "https://aiesec.org/opportunity/global-talent/1327149", "Data Engineering Intern", "Nugegoda, Sri Lanka", "No", "33 applicants", "3 - 6 Months", "Altria Consulting (PVT) LTD")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    # Column A holds a numeric-looking opportunity id; force text storage
    # (NumberFormat "@" before the write) then drop back to the default
    # "Normal" style so no stray number-format style is left on the cell.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}
